$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1: view state only (scrolled / selection moved down)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 66
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("B80").Select()

# ---------------------------------------------------------------------
# Sheet2: move the "OK?" row down, add new rows / column
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Capture the existing "OK?" row content before relocating it
$okLabel   = $ws2.Range("A13").Value2
$okFormula = $ws2.Range("B13").Formula
$ws2.Range("A13").ClearContents()
$ws2.Range("B13").ClearContents()

# New shared-string labels must be introduced in this order so the
# sharedStrings table comes out in the same sequence as the target file:
# "Super Tax Amount", "L's Total", "L Filter", "TEST"
$ws2.Range("A15").Value = "Super Tax Amount"
$ws2.Range("A17").Value = "L's Total"
$ws2.Range("A12").Value = "L Filter"
$ws2.Range("A16").Value = "TEST"

# Row 12: "L Filter" / "L1"
$ws2.Range("B12").Value = "L1"

# Row 14: relocated "OK?" row
$ws2.Range("A14").Value = $okLabel
$ws2.Range("B14").Formula = $okFormula

# Row 15: Super Tax Amount
$ws2.Range("B15").Formula = '=SUMIF(J4:J6,">200")'

# Row 16: TEST
$ws2.Range("B16").Formula = '=SUMIF(J4:J6,">" & J4)'

# Row 17: L's Total
$ws2.Range("B17").Formula = "=SUMIF(E4:E6,B12,J4:J6)"

# Row 18: plain filtered sum (no label)
$ws2.Range("B18").Formula = '=SUMIF(E4:E6,"L*",J4:J6)'

# Widen column A to fit the new labels
$ws2.Columns.Item(1).ColumnWidth = 15.706730769230772

$ws2.Range("B16").Select()

# ---------------------------------------------------------------------
# Sheet3: wipe the scratch/date-math content, rewrite as a small SUMIF demo
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Clear()

$ws3.Range("A1").Value = 100
$ws3.Range("A2").Value = 200
$ws3.Range("A3").Value = 300
$ws3.Range("B3").Formula = '=SUMIF(A1:A3,">"&A1)'

$ws3.Range("B3").Select()

$wb.Save()
